# Applies the diff: insert a new data row at sheet row 9 (shifting existing
# rows 9-49 down to 10-50) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; this shifts rows 9..49 down to 10..50
# and copies formatting (incl. the date style in column D) from the row above.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's values.
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44602
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100108
$ws.Cells.Item(9, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(9, 9).Value = 100108004
$ws.Cells.Item(9, 10).Value = "Papaya"
$ws.Cells.Item(9, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 60
$ws.Cells.Item(9, 14).Value = 23000
$ws.Cells.Item(9, 15).Value = 23000
$ws.Cells.Item(9, 16).Value = 23000
$ws.Cells.Item(9, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 19).Value = 2300
$ws.Cells.Item(9, 20).Value = 10
